$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Generate Report for Handback
#
# For each language sheet (zh-cn, de-de) and each data row, fill in the
# "Latest Target File" (F) and "Latest Handback File" (G) columns - mirroring
# the Source File Name (A) and Latest Handoff File (D) values/hyperlinks -
# and stamp the "Latest Handback DateTime" (H) column with the real handback
# timestamp. Also flip the Status text everywhere it says "Ready for
# handoff" to "Handed back: in sync with en-US".
# ---------------------------------------------------------------------------

$oldStatus = "Ready for handoff"
$newStatus = "Handed back: in sync with en-US"

# Handback timestamps per language sheet (row 2 and row 3 share the same
# handback run time on each sheet).
$handbackTimes = @{
    "zh-cn" = "2016-03-13 15:03:11"
    "de-de" = "2016-03-13 15:03:18"
}

function Set-StatusColumn($ws, $col) {
    for ($r = 2; $r -le 3; $r++) {
        $cell = $ws.Cells.Item($r, $col)
        if ($cell.Value2 -eq $oldStatus) {
            $cell.Value = $newStatus
        }
    }
}

# --- Overview sheet: just flip the Status text in columns B (zh-cn) and C (de-de) ---
$overview = $wb.Worksheets.Item("Overview")
Set-StatusColumn $overview 2
Set-StatusColumn $overview 3

foreach ($langName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($langName)

    # Status column (C) on the language sheet itself.
    Set-StatusColumn $ws 3

    for ($r = 2; $r -le 3; $r++) {
        $srcCell = $ws.Cells.Item($r, 1)   # A: Source File Name
        $handoffCell = $ws.Cells.Item($r, 4) # D: Latest Handoff File
        $targetCell = $ws.Cells.Item($r, 6)  # F: Latest Target File
        $handbackCell = $ws.Cells.Item($r, 7) # G: Latest Handback File
        $handbackDateCell = $ws.Cells.Item($r, 8) # H: Latest Handback DateTime

        $srcValue = $srcCell.Value2
        $handoffValue = $handoffCell.Value2

        # Find the existing hyperlinks on A<r> and D<r> so the new F/G
        # hyperlinks can point at the same targets/display text.
        $srcHyperlink = $null
        $handoffHyperlink = $null
        foreach ($hl in $ws.Hyperlinks) {
            $addr = $hl.Range.Address()
            if ($addr -eq ("$" + "A" + "$" + $r)) { $srcHyperlink = $hl }
            if ($addr -eq ("$" + "D" + "$" + $r)) { $handoffHyperlink = $hl }
        }

        # F<r>: Latest Target File == Source File Name
        $targetCell.Value = $srcValue
        $targetCell.Style = $srcCell.Style
        if ($srcHyperlink -ne $null) {
            $ws.Hyperlinks.Add($targetCell, $srcHyperlink.Address, "", "", $srcHyperlink.TextToDisplay) | Out-Null
        }

        # G<r>: Latest Handback File == Latest Handoff File
        $handbackCell.Value = $handoffValue
        $handbackCell.Style = $handoffCell.Style
        if ($handoffHyperlink -ne $null) {
            $ws.Hyperlinks.Add($handbackCell, $handoffHyperlink.Address, "", "", $handoffHyperlink.TextToDisplay) | Out-Null
        }

        # H<r>: Latest Handback DateTime -> real handback timestamp
        $handbackDateCell.Value = $handbackTimes[$langName]
    }
}

Write-Output "Generated handback report"
